$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BB1: copy the date-header format from BA1 (xlPasteFormats), then set its value
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("BB1").Value = 45986

# BB3:BB18 - duplicate the value already present in the corresponding BA cell
$ws.Range("BB3").Value = -5.109987415979145
$ws.Range("BB4").Value = 2.253603114136604
$ws.Range("BB5").Value = 3.738382206110891
$ws.Range("BB6").Value = 1.165974434765671
$ws.Range("BB7").Value = -0.0426719751787874
$ws.Range("BB8").Value = 1.529758493743438
$ws.Range("BB9").Value = 1.358758534900462
$ws.Range("BB10").Value = 1.664905435092301
$ws.Range("BB11").Value = 2.145670176886982
$ws.Range("BB12").Value = 1.976124254426503
$ws.Range("BB13").Value = 0.7060158009350337
$ws.Range("BB14").Value = -4.180878843351332
$ws.Range("BB15").Value = 1.312484974417294
$ws.Range("BB16").Value = 2.386394320099283
$ws.Range("BB17").Value = 0.2104414886460626
$ws.Range("BB18").Value = -0.3095793941792935

# BB19:BB21 - new forecast values (distinct from BA column)
$ws.Range("BB19").Value = -0.08656168856399082
$ws.Range("BB20").Value = -0.1516437243033186
$ws.Range("BB21").Value = -0.1967561196116963

Write-Host "done"
